$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 52.08999
$ws.Range("H2").Value = 156.26997
$ws.Range("I2").Value = 0.9401105828221099
$ws.Range("J2").Value = 0.9401105828221098
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.06447966666666667
$ws.Range("N2").Value = 0.193439
$ws.Range("O2").Value = 0.001101138907643723
$ws.Range("P2").Value = 0.001101138907643722
$ws.Range("Q2").Value = 3.35874519187
$ws.Range("R2").Value = 30.22870672683
$ws.Range("S2").Value = 0.001035192340233041
$ws.Range("T2").Value = 0.001035192340233041
$ws.Range("G3").Value = 52.08999
$ws.Range("H3").Value = 156.26997
$ws.Range("I3").Value = 0.9401105828221099
$ws.Range("J3").Value = 0.9401105828221098
$ws.Range("O3").Value = 0.00657695954769643
$ws.Range("P3").Value = 0.006576959547696431
$ws.Range("Q3").Value = 20.06134839538
$ws.Range("R3").Value = 180.55213555842
$ws.Range("S3").Value = 0.006183069273582332
$ws.Range("T3").Value = 0.006183069273582332
$ws.Range("G4").Value = 52.08999
$ws.Range("H4").Value = 156.26997
$ws.Range("I4").Value = 0.9401105828221099
$ws.Range("J4").Value = 0.9401105828221098
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.008175
$ws.Range("N4").Value = 0.024525
$ws.Range("O4").Value = 0.0001396069650378791
$ws.Range("P4").Value = 0.0001396069650378791
$ws.Range("Q4").Value = 0.42583566825
$ws.Range("R4").Value = 3.83252101425
$ws.Range("S4").Value = 0.0001312459852677864
$ws.Range("T4").Value = 0.0001312459852677864
$ws.Range("G5").Value = 52.08999
$ws.Range("H5").Value = 156.26997
$ws.Range("I5").Value = 0.9401105828221099
$ws.Range("J5").Value = 0.9401105828221098
$ws.Range("M5").Value = 58.099467
$ws.Range("N5").Value = 174.298401
$ws.Range("O5").Value = 0.992182294579622
$ws.Range("P5").Value = 0.992182294579622
$ws.Range("Q5").Value = 3026.40065503533
$ws.Range("R5").Value = 27237.60589531797
$ws.Range("S5").Value = 0.9327610752230268
$ws.Range("T5").Value = 0.9327610752230268
$ws.Range("I6").Value = 0.009851545038079508
$ws.Range("J6").Value = 0.009851545038079508
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.06447966666666667
$ws.Range("N6").Value = 0.193439
$ws.Range("O6").Value = 0.001101138907643723
$ws.Range("P6").Value = 0.001101138907643722
$ws.Range("Q6").Value = 0.03519674188733334
$ws.Range("R6").Value = 0.316770676986
$ws.Range("S6").Value = [double]"1.084791954183381E-05"
$ws.Range("T6").Value = [double]"1.08479195418338E-05"
$ws.Range("I7").Value = 0.009851545038079508
$ws.Range("J7").Value = 0.009851545038079508
$ws.Range("O7").Value = 0.00657695954769643
$ws.Range("P7").Value = 0.006576959547696431
$ws.Range("S7").Value = [double]"6.479321319775841E-05"
$ws.Range("T7").Value = [double]"6.479321319775843E-05"
$ws.Range("I8").Value = 0.009851545038079508
$ws.Range("J8").Value = 0.009851545038079508
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.008175
$ws.Range("N8").Value = 0.024525
$ws.Range("O8").Value = 0.0001396069650378791
$ws.Range("P8").Value = 0.0001396069650378791
$ws.Range("Q8").Value = 0.00446238915
$ws.Range("R8").Value = 0.04016150235
$ws.Range("S8").Value = [double]"1.375344303700257E-06"
$ws.Range("T8").Value = [double]"1.375344303700257E-06"
$ws.Range("I9").Value = 0.009851545038079508
$ws.Range("J9").Value = 0.009851545038079508
$ws.Range("M9").Value = 58.099467
$ws.Range("N9").Value = 174.298401
$ws.Range("O9").Value = 0.992182294579622
$ws.Range("P9").Value = 0.992182294579622
$ws.Range("Q9").Value = 31.714058857686
$ws.Range("R9").Value = 285.426529719174
$ws.Range("S9").Value = 0.009774528561036217
$ws.Range("T9").Value = 0.009774528561036217
$ws.Range("G10").Value = 1.744358333333333
$ws.Range("H10").Value = 5.233075
$ws.Range("I10").Value = 0.03148185917103467
$ws.Range("J10").Value = 0.03148185917103467
$ws.Range("K10").Value = 1
$ws.Range("L10").Value = 0.3333333333333333
$ws.Range("M10").Value = 0.06447966666666667
$ws.Range("N10").Value = 0.193439
$ws.Range("O10").Value = 0.001101138907643723
$ws.Range("P10").Value = 0.001101138907643722
$ws.Range("Q10").Value = 0.1124756438805556
$ws.Range("R10").Value = 1.012280794925
$ws.Range("S10").Value = [double]"3.466590001818663E-05"
$ws.Range("T10").Value = [double]"3.466590001818661E-05"
$ws.Range("G11").Value = 1.744358333333333
$ws.Range("H11").Value = 5.233075
$ws.Range("I11").Value = 0.03148185917103467
$ws.Range("J11").Value = 0.03148185917103467
$ws.Range("O11").Value = 0.00657695954769643
$ws.Range("P11").Value = 0.006576959547696431
$ws.Range("Q11").Value = 0.6718023991055556
$ws.Range("R11").Value = 6.04622159195
$ws.Range("S11").Value = 0.0002070549142541709
$ws.Range("T11").Value = 0.0002070549142541709
$ws.Range("G12").Value = 1.744358333333333
$ws.Range("H12").Value = 5.233075
$ws.Range("I12").Value = 0.03148185917103467
$ws.Range("J12").Value = 0.03148185917103467
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.008175
$ws.Range("N12").Value = 0.024525
$ws.Range("O12").Value = 0.0001396069650378791
$ws.Range("P12").Value = 0.0001396069650378791
$ws.Range("Q12").Value = 0.014260129375
$ws.Range("R12").Value = 0.128341164375
$ws.Range("S12").Value = [double]"4.39508681261807E-06"
$ws.Range("T12").Value = [double]"4.395086812618069E-06"
$ws.Range("G13").Value = 1.744358333333333
$ws.Range("H13").Value = 5.233075
$ws.Range("I13").Value = 0.03148185917103467
$ws.Range("J13").Value = 0.03148185917103467
$ws.Range("M13").Value = 58.099467
$ws.Range("N13").Value = 174.298401
$ws.Range("O13").Value = 0.992182294579622
$ws.Range("P13").Value = 0.992182294579622
$ws.Range("Q13").Value = 101.346289423675
$ws.Range("R13").Value = 912.1166048130751
$ws.Range("S13").Value = 0.0312357432699497
$ws.Range("T13").Value = 0.03123574326994969
$ws.Range("G14").Value = 1.028158333333334
$ws.Range("H14").Value = 3.084475
$ws.Range("I14").Value = 0.01855601296877595
$ws.Range("J14").Value = 0.01855601296877594
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 0.3333333333333333
$ws.Range("M14").Value = 0.06447966666666667
$ws.Range("N14").Value = 0.193439
$ws.Range("O14").Value = 0.001101138907643723
$ws.Range("P14").Value = 0.001101138907643722
$ws.Range("Q14").Value = 0.06629530661388891
$ws.Range("R14").Value = 0.5966577595250001
$ws.Range("S14").Value = [double]"2.043274785066069E-05"
$ws.Range("T14").Value = [double]"2.043274785066068E-05"
$ws.Range("G15").Value = 1.028158333333334
$ws.Range("H15").Value = 3.084475
$ws.Range("I15").Value = 0.01855601296877595
$ws.Range("J15").Value = 0.01855601296877594
$ws.Range("O15").Value = 0.00657695954769643
$ws.Range("P15").Value = 0.006576959547696431
$ws.Range("Q15").Value = 0.3959732480388889
$ws.Range("R15").Value = 3.56375923235
$ws.Range("S15").Value = 0.0001220421466621697
$ws.Range("T15").Value = 0.0001220421466621697
$ws.Range("G16").Value = 1.028158333333334
$ws.Range("H16").Value = 3.084475
$ws.Range("I16").Value = 0.01855601296877595
$ws.Range("J16").Value = 0.01855601296877594
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.008175
$ws.Range("N16").Value = 0.024525
$ws.Range("O16").Value = 0.0001396069650378791
$ws.Range("P16").Value = 0.0001396069650378791
$ws.Range("Q16").Value = 0.008405194375000001
$ws.Range("R16").Value = 0.07564674937500002
$ws.Range("S16").Value = [double]"2.590548653774334E-06"
$ws.Range("T16").Value = [double]"2.590548653774334E-06"
$ws.Range("G17").Value = 1.028158333333334
$ws.Range("H17").Value = 3.084475
$ws.Range("I17").Value = 0.01855601296877595
$ws.Range("J17").Value = 0.01855601296877594
$ws.Range("M17").Value = 58.099467
$ws.Range("N17").Value = 174.298401
$ws.Range("O17").Value = 0.992182294579622
$ws.Range("P17").Value = 0.992182294579622
$ws.Range("Q17").Value = 59.73545115827501
$ws.Range("R17").Value = 537.619060424475
$ws.Range("S17").Value = 0.01841094752560934
$ws.Range("T17").Value = 0.01841094752560934
